$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("League Summary")

# Insert a new column before column A, shifting the existing "Team" /
# "Games Played" / ... columns (and their header/value cells) one
# position to the right.
$ws.Columns("A:A").Insert()

# Populate the new column A with the "League" header and "Fray" value.
$ws.Range("A1").Value = "League"
$ws.Range("A2").Value = "Fray"

# Give the new header cell (A1) the same formatting (bold font, border,
# center/top alignment) already used by the rest of row 1 -- copy it from
# the neighboring header cell so the same style record is reused instead
# of synthesizing a near-duplicate one.
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)

# Column widths after the insert: new col A = 20 (matches the old "Team"
# column that shifted into B), and every other column keeps its original
# width, shifted one slot to the right; the new trailing col H = 12.
# NOTE: ColumnWidth applies Excel's char-width -> internal-unit padding
# (adds 5/6) on write in this runtime, so the input is pre-compensated
# by subtracting 5/6 to land exactly on the intended integer widths.
$ws.Columns("A:A").ColumnWidth = 19.166666666666668  # -> 20
$ws.Columns("B:B").ColumnWidth = 19.166666666666668  # -> 20
$ws.Columns("C:C").ColumnWidth = 14.166666666666666  # -> 15
$ws.Columns("D:D").ColumnWidth = 15.166666666666666  # -> 16
$ws.Columns("E:E").ColumnWidth = 11.166666666666666  # -> 12
$ws.Columns("F:F").ColumnWidth = 11.166666666666666  # -> 12
$ws.Columns("G:G").ColumnWidth = 11.166666666666666  # -> 12
$ws.Columns("H:H").ColumnWidth = 11.166666666666666  # -> 12
